# margin notes working in OTA docx generated
#
# Switches the document's base font (and a couple of heading-linked
# styles) over to Baskerville, and adds a new "marginOuter" /
# "marginOuterChar" style pair used for margin notes.

$d = $word.ActiveDocument
$styles = $d.Styles

# --- Normal: literal Baskerville ascii/hAnsi font -----------------------
$normal = $styles.Item("Normal")
$normal.Font.NameAscii = "Baskerville"
$normal.Font.NameOther = "Baskerville"

# --- Heading1 (paragraph style linked to Heading1Char) -------------------
# The authoritative edit drops the ascii/hAnsi *theme* font reference so
# the heading now inherits the literal Baskerville family from Normal.
$heading1 = $styles.Item("Heading1")
$heading1.Font.NameAscii = "Baskerville"
$heading1.Font.NameOther = "Baskerville"

# --- Heading1Char (linked character style) --------------------------------
# Picks up literal Baskerville ascii/hAnsi while keeping the east-asian /
# complex-script theme references untouched.
$heading1Char = $styles.Item("Heading1Char")
$heading1Char.Font.NameAscii = "Baskerville"
$heading1Char.Font.NameOther = "Baskerville"

# --- GeneratedTitle: literal Baskerville ascii/hAnsi font -----------------
$generatedTitle = $styles.Item("GeneratedTitle")
$generatedTitle.Font.NameAscii = "Baskerville"
$generatedTitle.Font.NameOther = "Baskerville"

# --- New styles: marginOuter / marginOuterChar ----------------------------
$marginOuter = $styles.Add("marginOuter", 1)
$marginOuter.BaseStyle = "Normal"
$marginOuter.NextParagraphStyle = "Normal"
$marginOuter.QuickStyle = $true
$marginOuter.Font.Size = 10
$marginOuter.Font.SizeBi = 10

$marginOuterChar = $styles.Add("marginOuterChar", 2)
$marginOuterChar.NameLocal = "marginOuter Char"
$marginOuterChar.BaseStyle = "DefaultParagraphFont"
$marginOuterChar.Font.NameAscii = "Baskerville"
$marginOuterChar.Font.NameOther = "Baskerville"
$marginOuterChar.Font.Size = 10
$marginOuterChar.Font.SizeBi = 10

# Link the paragraph <-> character style pair both ways.
$marginOuter.LinkStyle = "marginOuterChar"
$marginOuterChar.LinkStyle = "marginOuter"
